# Weekly update: a new price record for "Apio" (Macroferia Regional de Talca)
# is inserted as the new row 32. All existing records from row 32 downward
# are pushed down by one row (handled automatically by the row Insert),
# and the table grows from A1:R149 to A1:R150.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 32 - shifts old rows 32..149 to 33..150
$ws.Rows("32:32").Insert()

# Populate the newly inserted row 32 with the new record
$ws.Range("A32").Value = 5
$ws.Range("B32").Value = "Macroferia Regional de Talca"
$ws.Range("C32").Value = "Maule"
$ws.Range("D32").Value = 44659
$ws.Range("E32").Value = 7
$ws.Range("F32").Value = 100112017
$ws.Range("G32").Value = "Apio"
$ws.Range("H32").Value = "Americana (o)"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 500
$ws.Range("K32").Value = 7000
$ws.Range("L32").Value = 7000
$ws.Range("M32").Value = 7000
$ws.Range("N32").Value = '$/docena de matas'
$ws.Range("O32").Value = "Provincia del Elquí"
$ws.Range("P32").Value = 1167
$ws.Range("Q32").Value = 6
$ws.Range("R32").Value = "Hortaliza"
